$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Projections")

# Update row 4 values (FS and KZN updates)
$ws.Range("A4").Value = 37530
$ws.Range("B4").Value = 845
$ws.Range("C4").Value = 97042
$ws.Range("D4").Value = 0.8913436097821972
$ws.Range("E4").Value = 0.02007217454637576
$ws.Range("F4").Value = 2.304756198585654
$ws.Range("G4").Value = 0.2233735931670334
$ws.Range("H4").Value = 0.01324390892681846
$ws.Range("I4").Value = 0.3589482933770215
$ws.Range("J4").Value = 0.02860994820120611
$ws.Range("K4").Value = 0.009342840495265839
$ws.Range("L4").Value = 0.05376641014642117
$ws.Range("M4").Value = 323153
$ws.Range("N4").Value = 263
$ws.Range("O4").Value = 682982
$ws.Range("P4").Value = 8.023739693690269
$ws.Range("Q4").Value = 0.006552486261007747
$ws.Range("R4").Value = 16.95812267353212
$ws.Range("S4").Value = 0.5169970601880156
$ws.Range("T4").Value = 0.004890021849963582
$ws.Range("U4").Value = 0.6714657200384131
$ws.Range("V4").Value = 0.2227171165814395
$ws.Range("W4").Value = 0.009465147851420249
$ws.Range("X4").Value = 0.4669746120339427
